# Update strain names and experiment conditions in bioSample sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row groups (startRow, endRow, genotype, strain)
$groups = @(
    @(2,  4,  "CNAG_00000",            "KN99alpha"),
    @(5,  7,  "CNAG_02566",            "TDY1452"),
    @(8,  10, "CNAG_05222",            "TDY1118"),
    @(11, 13, "CNAG_01438",            "TDY1451"),
    @(14, 16, "CNAG_05222.CNAG_02566", "TDY1665"),
    @(17, 19, "CNAG_05222.CNAG_01438", "TDY1652"),
    @(20, 22, "CNAG_02566.CNAG_01438", "TDY1665")
)

foreach ($g in $groups) {
    $startRow = $g[0]
    $endRow = $g[1]
    $genotype = $g[2]
    $strain = $g[3]
    for ($r = $startRow; $r -le $endRow; $r++) {
        $ws.Cells.Item($r, 4).Value = "90minuteInduction"  # column D: experimentDesign
        $ws.Cells.Item($r, 2).Value = "S.GISH"             # column B: harvester
        $ws.Cells.Item($r, 6).Value = $strain              # column F: strain
        $ws.Cells.Item($r, 7).Value = $genotype            # column G: genotype
    }
}

# Restore the view's selected range to match the saved workbook state.
$null = $ws.Range("F21:F22").Select()
